$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.04416460176223325
$ws.Range("E2").Value = 0.4396191691991618
$ws.Range("F2").Value = 1.528290297222981
$ws.Range("G2").Value = 1.673767081463325
$ws.Range("H2").Value = 1.010542361786406
$ws.Range("D3").Value = 0.03830832642181292
$ws.Range("E3").Value = 0.382347927903524
$ws.Range("F3").Value = 1.363286054659937
$ws.Range("G3").Value = 1.472460773835564
$ws.Range("H3").Value = 0.9231199012414208
$ws.Range("D4").Value = 0.0347263708503931
$ws.Range("E4").Value = 0.3474053453013681
$ws.Range("F4").Value = 1.262945067803372
$ws.Range("G4").Value = 1.349747030403478
$ws.Range("H4").Value = 0.8702008080870769
$ws.Range("D5").Value = 0.03326964694163337
$ws.Range("E5").Value = 0.3332162809223007
$ws.Range("F5").Value = 1.222288067622259
$ws.Range("G5").Value = 1.299949458491483
$ws.Range("H5").Value = 0.8488200623750686
$ws.Range("D6").Value = 0.03302792226328677
$ws.Range("E6").Value = 0.3308630769826237
$ws.Range("F6").Value = 1.215550724464322
$ws.Range("G6").Value = 1.291692836908965
$ws.Range("H6").Value = 0.8452807431613678
$ws.Range("D7").Value = 0.03470671369254319
$ws.Range("E7").Value = 0.34721379011431
$ws.Range("F7").Value = 1.262395826126095
$ws.Range("G7").Value = 1.349074614279175
$ws.Range("H7").Value = 0.8699117225502562
$ws.Range("D8").Value = 0.04214216942422411
$ws.Range("E8").Value = 0.4198225367166231
$ws.Range("F8").Value = 1.471187856943203
$ws.Range("G8").Value = 1.604163365212287
$ws.Range("H8").Value = 0.9802376235733732
$ws.Range("D9").Value = 0.05685744970708129
$ws.Range("E9").Value = 0.5642327170640442
$ws.Range("F9").Value = 1.888902018709416
$ws.Range("G9").Value = 2.112127782647519
$ws.Range("H9").Value = 1.202918706961327
$ws.Range("D10").Value = 0.06778656051734799
$ws.Range("E10").Value = 0.671944623436147
$ws.Range("F10").Value = 2.201659492699918
$ws.Range("G10").Value = 2.491042104895712
$ws.Range("H10").Value = 1.370841485233484
$ws.Range("D11").Value = 0.07279217338034982
$ws.Range("E11").Value = 0.7213816670721229
$ws.Range("F11").Value = 2.345399064368479
$ws.Range("G11").Value = 2.664885521213535
$ws.Range("H11").Value = 1.44827547428082
$ws.Range("D12").Value = 0.07469321267008411
$ws.Range("E12").Value = 0.7401723167000966
$ws.Range("F12").Value = 2.400055170825851
$ws.Range("G12").Value = 2.730945692542946
$ws.Range("H12").Value = 1.477756451588164
$ws.Range("D13").Value = 0.07428353371486196
$ws.Range("E13").Value = 0.7361221926845047
$ws.Range("F13").Value = 2.388273756458318
$ws.Range("G13").Value = 2.716707962544604
$ws.Range("H13").Value = 1.471400018306952
$ws.Range("D14").Value = 0.07294845857009591
$ws.Range("E14").Value = 0.7229261437162222
$ws.Range("F14").Value = 2.349891060477688
$ws.Range("G14").Value = 2.670315630292464
$ws.Range("H14").Value = 1.450697669861938
$ws.Range("D15").Value = 0.07213142521979421
$ws.Range("E15").Value = 0.7148524956227789
$ws.Range("F15").Value = 2.326410295331868
$ws.Range("G15").Value = 2.641929398215893
$ws.Range("H15").Value = 1.438037758433495
$ws.Range("D16").Value = 0.06746016902407348
$ws.Range("E16").Value = 0.6687232001522858
$ws.Range("F16").Value = 2.19229646974
$ws.Range("G16").Value = 2.479712166869149
$ws.Range("H16").Value = 1.365802710240814
$ws.Range("D17").Value = 0.06460364233217319
$ws.Range("E17").Value = 0.6405414693057452
$ws.Range("F17").Value = 2.110407240501445
$ws.Range("G17").Value = 2.38058685249996
$ws.Range("H17").Value = 1.321762217021217
$ws.Range("D18").Value = 0.06296379515427475
$ws.Range("E18").Value = 0.6243729184996312
$ws.Range("F18").Value = 2.063443689084238
$ws.Range("G18").Value = 2.323710192705562
$ws.Range("H18").Value = 1.2965291561203
$ws.Range("D19").Value = 0.06240909297230246
$ws.Range("E19").Value = 0.6189053326470031
$ws.Range("F19").Value = 2.047565729167587
$ws.Range("G19").Value = 2.30447583926275
$ws.Range("H19").Value = 1.288002240429705
$ws.Range("D20").Value = 0.06490739369478149
$ws.Range("E20").Value = 0.6435371872360776
$ws.Range("F20").Value = 2.11911020998528
$ws.Range("G20").Value = 2.391124544869797
$ws.Range("H20").Value = 1.326440215915341
$ws.Range("D21").Value = 0.07334044737078216
$ws.Range("E21").Value = 0.7268001897907368
$ws.Range("F21").Value = 2.361158763589515
$ws.Range("G21").Value = 2.683935813114488
$ws.Range("H21").Value = 1.45677408627796
$ws.Range("D22").Value = 0.07888443913446963
$ws.Range("E22").Value = 0.7816282444840681
$ws.Range("F22").Value = 2.520670277776929
$ws.Range("G22").Value = 2.876650817821201
$ws.Range("H22").Value = 1.542882076632281
$ws.Range("D23").Value = 0.07592231310700015
$ws.Range("E23").Value = 0.7523255529325468
$ws.Range("F23").Value = 2.435410421870188
$ws.Range("G23").Value = 2.773666093207794
$ws.Range("H23").Value = 1.496837023020134
$ws.Range("D24").Value = 0.06477006025221499
$ws.Range("E24").Value = 0.6421827193270531
$ws.Range("F24").Value = 2.115175238321655
$ws.Range("G24").Value = 2.386360107227858
$ws.Range("H24").Value = 1.324325024093469
$ws.Range("D25").Value = 0.05285854641977039
$ws.Range("E25").Value = 0.5249106469200768
$ws.Range("F25").Value = 1.774930254878257
$ws.Range("G25").Value = 1.973780395381766
$ws.Range("H25").Value = 1.141953052602048
